$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be force-typed as text
# (matching the source data, which stores all Price/Volume cells as strings)
# by toggling NumberFormat to Text ("@") around the assignment, then restoring it.

$ws.Range("D2").Value = '43.830.07'
$ws.Range("D3").Value = '2.295.09'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.21'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +17.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.24'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.54'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +7.37%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.03'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +14.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.82'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '2.639.32'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.855'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("D17").Value = '2.293.81'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = '43.736.06'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.13'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +14.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.37'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.85'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +7.71%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.16'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("E25").Value = '  +7.46%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.67'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +2.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '43.13'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +12.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.90'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '175.15'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("E33").Value = '  +4.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.61'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.69'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +3.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.83'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +2.96%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0362'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +2.59%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.106'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -3.19%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.84'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +8.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.41'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +18.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '75.98'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +17.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.242'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +2.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.34'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +21.29%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.40'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +1.88%  '
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.26'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +2.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.22'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +3.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0995'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -2.78%  '
